$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (s="2", used by A2:A464) down onto the new rows
# A465:A491 before writing values, so the new date cells keep the same
# centered/bordered "YYYY-MM-DD HH:MM:SS" look as the rest of column A.
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)

# New daily rows (date serial, nuovi pos., somma mobile 7gg., per 100mila ab.)
# continuing the series through 2022-01-06 ("aggiornamento fino a 6 gennaio 2022").
$newRows = @(
    @(465, 44539, 4, 9, 96.9409737182249),
    @(466, 44540, 1, 10, 107.7121930202499),
    @(467, 44541, 0, 8, 86.16975441619991),
    @(468, 44542, 6, 13, 140.0258509263249),
    @(469, 44543, 6, 18, 193.8819474364498),
    @(470, 44544, 0, 18, 193.8819474364498),
    @(471, 44545, 0, 17, 183.1107281344248),
    @(472, 44546, 19, 32, 344.6790176647996),
    @(473, 44547, 16, 47, 506.2473071951745),
    @(474, 44548, 5, 52, 560.1034037052995),
    @(475, 44550, 6, 52, 560.1034037052995),
    @(476, 44551, 5, 51, 549.3321844032745),
    @(477, 44552, 1, 52, 560.1034037052995),
    @(478, 44553, 9, 61, 657.0443774235243),
    @(479, 44554, 7, 49, 527.7897457992245),
    @(480, 44555, 6, 39, 420.0775527789746),
    @(481, 44556, 10, 44, 473.9336492890995),
    @(482, 44557, 12, 50, 538.5609651012494),
    @(483, 44558, 0, 45, 484.7048685911245),
    @(484, 44559, 3, 47, 506.2473071951745),
    @(485, 44560, 13, 51, 549.3321844032745),
    @(486, 44561, 18, 62, 667.8155967255493),
    @(487, 44562, 15, 71, 764.7565704437743),
    @(488, 44563, 7, 68, 732.4429125376993),
    @(489, 44564, 2, 58, 624.7307195174493),
    @(490, 44565, 3, 61, 657.0443774235243),
    @(491, 44566, 6, 64, 689.3580353295993)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
